$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 38 (a new prospect, "NAI LEGACY"), shifting the
# existing rows 38:41 down to 39:42.
$ws.Rows.Item(38).Insert()

$ws.Cells.Item(38, 1).Value = "NAI LEGACY"
$ws.Cells.Item(38, 2).Value = "Steiner, Owen A"
$ws.Cells.Item(38, 3).Value = "015"
$ws.Cells.Item(38, 5).Value = "0008398"

# Match the row height used throughout the rest of the table.
$ws.Rows.Item(38).RowHeight = 13.05

# Column D (Last Invoice Date) stays blank for a prospect with no invoice
# yet; column F is also blank on every row. Copy the empty, unstyled
# F-cell from the row above so the new row keeps the same blank-cell shape
# as its neighbours instead of simply having no F cell at all.
$ws.Cells.Item(37, 6).Copy($ws.Cells.Item(38, 6))
